# Defect Tracker template update:
#   - The defect tracker's assignee name was corrected from "Arjun" to
#     "Aryan" throughout the "Assigned To" column (F2:F9) on the
#     "Defects" sheet.
#   - The active selection is left on F10 (the first empty cell below the
#     "Assigned To" data), matching where the editor's cursor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Defects")

# Correct the assignee's name for every defect row.
$ws.Range("F2:F9").Value = "Aryan"

# Leave the selection where the edit session ended.
$ws.Range("F10").Select()
